# Add a "2022-Q1" sheet (before "总计") with one fund holding row, and
# refresh the "总计" (summary) sheet so it includes the new quarter.
#
# Sheet-id bookkeeping note: new sheets are assigned sheetId = (max of
# current sheetIds) + 1. To land on the same numbering as the target
# workbook (2022-Q1 -> sheetId 4, 总计 -> sheetId 5) we:
#   1. add a brand-new placeholder sheet right after the existing "总计"
#      sheet (while "总计" still holds id 4, so the placeholder claims id 5
#      and is already positioned last, which is where the refreshed
#      "总计" sheet belongs);
#   2. repurpose the *existing* "总计" sheet object (which keeps its id 4
#      and is already positioned right where "2022-Q1" belongs) into the
#      new "2022-Q1" sheet;
#   3. rename the placeholder into the new "总计" sheet.

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# Step 1: placeholder for the refreshed 总计 sheet, placed right after the
# current 总计 (i.e. at the very end of the workbook).
$newTotal = $wb.Worksheets.Add($null, $oldTotal)
$newTotal.Name = "总计_new"

# Step 2: turn the old 总计 sheet into the new 2022-Q1 sheet.
$q1 = $oldTotal
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Bring over the header row + row-layout (styles) from an existing
# quarter sheet ("2021-Q4") which already carries the right header
# labels/format, then overwrite with 2022-Q1's own figures.
$q4.Range("A1:H2").Copy()
$q1.Range("A1:H2").PasteSpecial(-4122)   # xlPasteFormats
$q4.Range("A1:H2").Copy()
$q1.Range("A1:H2").PasteSpecial(-4163)   # xlPasteValues
$q1.Range("A1").Clear()                  # no A1 cell in the source layout

$q1.Range("A2").Value = 0

function Set-TextValue($ws, $addr, $text) {
    # Force a literal text value (even if it parses as a number/date),
    # matching the source workbook's inlineStr cells, without leaving a
    # stray "quote prefix" style behind.
    $ws.Range($addr).Formula = '="' + $text + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue $q1 "B2" "163302"
Set-TextValue $q1 "C2" "大摩资源优选混合(LOF)"
Set-TextValue $q1 "D2" "5.82"
Set-TextValue $q1 "E2" "81.78"
Set-TextValue $q1 "F2" "4.40"
Set-TextValue $q1 "G2" "0.2561"
$q1.Range("H2").Value = 3

# Step 3: populate the refreshed 总计 sheet.
$tot = $newTotal
$tot.Name = "总计"

$q4.Range("B1:D1").Copy()
$tot.Range("B1:D1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$tot.Range("A2:A5").PasteSpecial(-4122)

$tot.Range("B1").Value = "日期"
$tot.Range("C1").Value = "持有数量(只)"
$tot.Range("D1").Value = "持有市值(亿元)"

$tot.Range("A2").Value = 0
Set-TextValue $tot "B2" "2022-Q1"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 0.26

$tot.Range("A3").Value = 1
Set-TextValue $tot "B3" "2021-Q4"
$tot.Range("C3").Value = 3
$tot.Range("D3").Value = 0.48

$tot.Range("A4").Value = 2
Set-TextValue $tot "B4" "2021-Q1"
$tot.Range("C4").Value = 6
$tot.Range("D4").Value = 1.02

$tot.Range("A5").Value = 3
Set-TextValue $tot "B5" "2020-Q4"
$tot.Range("C5").Value = 3
$tot.Range("D5").Value = 0.54
